$wb = $excel.ActiveWorkbook

# Update the filter sheet's "Credit Cards" segmentation label.
# Previously this cell held the long description
# "Credit Cards: Silver, Gold, Platinum & Millennium"; it is now
# shortened to simply "Credit Cards" (DB validation driver change).
$wsFilter = $wb.Worksheets.Item("Cross_sell_Filter")
$wsFilter.Range("C2").Value = "Credit Cards"

# Reflect the user's last interaction: they ended up on the
# Cross_sell_Filter sheet with cell C12 selected (instead of the
# Cross_Sell_Report_EtoE sheet with B8 selected there).
$wsFilter.Activate()
$wsFilter.Range("C12").Select()
